$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2706
$ws1.Range("F4").Value = 1165
$ws1.Range("F5").Value = 1328
$ws1.Range("F6").Value = 289
$ws1.Range("F8").Value = 10469
$ws1.Range("F14").Value = 11959
$ws1.Range("F15").Value = 12364
$ws1.Range("F21").Value = 39

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2706
$ws4.Range("F5").Value = 1165
$ws4.Range("F6").Value = 1328
$ws4.Range("F7").Value = 289
$ws4.Range("F9").Value = 10469
$ws4.Range("F15").Value = 11959
$ws4.Range("F16").Value = 12364
$ws4.Range("F22").Value = 39
